$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix CN1 (row 4): Mid Y coordinate -6.5 -> -5, Rotation 180 -> 90
$ws.Range("C4").Value = -5
$ws.Range("E4").Value = 90

# Fix U1 (row 14): Rotation 90 -> 0
$ws.Range("E14").Value = 0

$ws.Range("E15").Select()
